$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in previously-empty spheroid measurements on row 14 (E14:H14) ---
# These four values are the ones that "fell off" the top of the I:L shift below.
$ws.Range("E14").Value = 96.7
$ws.Range("F14").Value = 65
$ws.Range("G14").Value = 60.9
$ws.Range("H14").Value = 768

# --- Shift multi_poly measurements (I24:L31) up by one row ---
# Row 24 <- old row 25, row 25 <- old row 26, ..., row 30 <- old row 31
# then row 31 is cleared.
$ws.Range("I24").Value = 46.4
$ws.Range("J24").Value = 31.8
$ws.Range("K24").Value = 31
$ws.Range("L24").Value = 64.400000000000006

$ws.Range("I25").Value = 83.4
$ws.Range("J25").Value = 67.8
$ws.Range("K25").Value = 59.6
$ws.Range("L25").Value = 442.5

$ws.Range("I26").Value = 58.6
$ws.Range("J26").Value = 54.6
$ws.Range("K26").Value = 51.8
$ws.Range("L26").Value = 210

$ws.Range("I27").Value = 97.3
$ws.Range("J27").Value = 65.599999999999994
$ws.Range("K27").Value = 54
$ws.Range("L27").Value = 594.1

$ws.Range("I28").Value = 69.099999999999994
$ws.Range("J28").Value = 48.9
$ws.Range("K28").Value = 45.9
$ws.Range("L28").Value = 223.7

$ws.Range("I29").Value = 91.1
$ws.Range("J29").Value = 63.3
$ws.Range("K29").Value = 45.8
$ws.Range("L29").Value = 371.5

$ws.Range("I30").Value = 47.2
$ws.Range("J30").Value = 36.6
$ws.Range("K30").Value = 34.700000000000003
$ws.Range("L30").Value = 113.8

$ws.Range("I31:L31").ClearContents()

# --- Update the view state: scroll position and active selection ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("F18").Select()
